# Update MSME Country Indicators - Czech Republic Summary figures to add
# a second decimal of precision. The source cells store these numbers as
# text (so they round-trip exactly, e.g. "85.39" rather than 85.39), so
# we force each cell's format to Text before writing the new value, then
# restore the cell's style to Normal (the sheet has no visible number
# formatting applied to these cells either way).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 13: Enterprises density (per 1000 people) -- Micro / SMEs / MSMEs
Set-TextValue $ws.Range("B13") "85.39"
Set-TextValue $ws.Range("C13") "3.91"
Set-TextValue $ws.Range("D13") "89.31"

# Row 14: Employment (% of total) -- Micro / SMEs / MSMEs
Set-TextValue $ws.Range("B14") "30.11"
Set-TextValue $ws.Range("C14") "38.36"
Set-TextValue $ws.Range("D14") "68.47"

# Row 16: Enterprises (% of total) -- Micro / SMEs / MSMEs
Set-TextValue $ws.Range("B16") "95.47"
Set-TextValue $ws.Range("C16") "4.38"
Set-TextValue $ws.Range("D16") "99.84"
